$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'246.61"
$ws.Range("D3").Formula = "'26.56"
$ws.Range("D4").Formula = "'5.094"
$ws.Range("D5").Formula = "'0.05612"
$ws.Range("D6").Formula = "'6.475"
$ws.Range("D7").Formula = "'0.8138"
$ws.Range("D8").Formula = "'0.8458"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Formula = "'0.009886"
$ws.Range("E9").Value = "8OneONEBestin24h"
$ws.Range("B10").Value = "BitrueCoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D10").Formula = "'0.02846"
$ws.Range("E10").Value = "9BitrueCoinBTR"
$ws.Range("B11").Value = "BitMartToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D11").Formula = "'0.09391"
$ws.Range("E11").Value = "10BitMartTokenBMX"
$ws.Range("B12").Value = "BitForexToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D12").Formula = "'0.001516"
$ws.Range("E12").Value = "11BitForexTokenBF"
$ws.Range("D13").Formula = "'0.006141"
$ws.Range("D17").Formula = "'0.3207"
$ws.Range("D19").Formula = "'0.06955"
$ws.Range("D20").Formula = "'0.03150"
$ws.Range("D21").Formula = "'0.1321"
$ws.Range("D22").Formula = "'3.747"
$ws.Range("D23").Formula = "'0.04642"
$ws.Range("D25").Formula = "'0.001248"
$ws.Range("D26").Formula = "'0.004591"
$ws.Range("E28").Value = "27UpBotsUBXT"
$ws.Range("D40").Formula = "'0.03668"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Formula = "'0.006233"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Formula = "'0.1056"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("D43").Formula = "'0.002605"
$ws.Range("D44").Formula = "'0.008912"
$ws.Range("D45").Formula = "'0.00005299"
$ws.Range("D48").Formula = "'0.002515"
